$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.375.88"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.844.03"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'240.33"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").Value = "'0.6338"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.07544"
$ws.Range("D9").Value = "'0.2957"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "'24.83"
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("D11").Value = "'0.07731"
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").Value = "'4.988"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "'0.6824"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "'83.02"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "'0.000009957"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").Value = "'6.160"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").Value = "29.396.88"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "'230.29"
$ws.Range("E18").Value = "  -3.64%  "
$ws.Range("D19").Value = "'12.44"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "'0.9996"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'7.542"
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  +233.90%  "
$ws.Range("D24").Value = "'16.62"
$ws.Range("E24").Value = "  +170.82%  "
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").Value = "'0.1400"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'8.372"
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("E29").Value = "  +172.57%  "
$ws.Range("D30").Value = "'1.467"
$ws.Range("D31").Value = "'0.05695"
$ws.Range("E31").Value = "  -3.98%  "
$ws.Range("D32").Value = "'1.252"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("D33").Value = "'4.119"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").Value = "'4.022"
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  -3.59%  "
$ws.Range("D36").Value = "'1.155"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("D37").Value = "'0.7142"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("D38").Value = "'2.600"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "1.244.20"
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("D40").Value = "'2.798"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("E42").Value = "  +266.10%  "
$ws.Range("D43").Value = "'0.9034"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").Value = "'0.9995"
$ws.Range("D45").Value = "'101.70"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").Value = "'66.26"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("D47").Value = "'7.057"
$ws.Range("E47").Value = "  -5.11%  "
$ws.Range("D48").Value = "'9.158"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").Value = "'0.4014"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D50").Value = "'1.700"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").Value = "'0.1122"
$ws.Range("E51").Value = "  -0.67%  "
